# Apply crypto price/volume updates to match the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.813.67'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.737.03'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''227.22'
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '''0.5138'
$ws.Range('E7').Value = '  +1.37%  '
$ws.Range('D8').Value = '''0.2684'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('D9').Value = '''39.27'
$ws.Range('E9').Value = '  -5.34%  '
$ws.Range('D10').Value = '''0.06081'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').Value = '1.734.67'
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').Value = '''0.06984'
$ws.Range('D13').Value = '''15.13'
$ws.Range('E13').Value = '  -3.13%  '
$ws.Range('D14').Value = '''0.6248'
$ws.Range('E14').Value = '  +3.68%  '
$ws.Range('D15').Value = '''4.487'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '''76.27'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = '''1.000'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').Value = '25.828.82'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '''11.41'
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('D21').Value = '''0.000006530'
$ws.Range('E21').Value = '  -4.55%  '
$ws.Range('D22').Value = '1.956.34'
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('D23').Value = '''4.022'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').Value = '''8.349'
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').Value = '''5.068'
$ws.Range('E25').Value = '  -2.57%  '
$ws.Range('D26').Value = '''136.49'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = '''1.505'
$ws.Range('E27').Value = '  +2.30%  '
$ws.Range('D28').Value = '''1.816'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = '''14.92'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('D30').Value = '''102.32'
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').Value = '''0.08286'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').Value = '''3.607'
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('D33').Value = '''3.348'
$ws.Range('E33').Value = '  -3.04%  '
$ws.Range('D34').Value = '''0.04404'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('D35').Value = '''2.610'
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('D36').Value = '''0.9705'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('D37').Value = '''0.5944'
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('D38').Value = '''2.680'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('D39').Value = '''0.01563'
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').Value = '''1.912'
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('D41').Value = '''0.9978'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').Value = '''101.72'
$ws.Range('E42').Value = '  -1.78%  '
$ws.Range('D43').Value = '''0.3787'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').Value = '''0.7281'
$ws.Range('E44').Value = '  -1.57%  '
$ws.Range('D45').Value = '''4.818'
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').Value = '''0.05476'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').Value = '''6.224'
$ws.Range('E47').Value = '  +4.20%  '
$ws.Range('D48').Value = '''0.1098'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('D49').Value = '''29.63'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('D50').Value = '''51.61'
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('E51').Value = '  +0.08%  '
